$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are stored as text so values like "1.002" or "11.20" are not
# reinterpreted as numbers (which would drop trailing zeros / change precision).
$priceCells = @("D2", "D3", "D4", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '23.315.20'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '1.628.11'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").Value = '297.78'
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("D7").Value = '0.3759'
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("D8").Value = '50.25'
$ws.Range("E8").Value = '  -2.05%  '
$ws.Range("D9").Value = '0.3471'
$ws.Range("E9").Value = '  -3.90%  '
$ws.Range("D10").Value = '0.08016'
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("E11").Value = '  -2.69%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '21.81'
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("D14").Value = '6.282'
$ws.Range("E14").Value = '  -2.94%  '
$ws.Range("D15").Value = '7.194'
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("D16").Value = '0.00001185'
$ws.Range("E16").Value = '  -3.60%  '
$ws.Range("D17").Value = '1.626.99'
$ws.Range("E17").Value = '  -1.38%  '
$ws.Range("D18").Value = '94.39'
$ws.Range("E18").Value = '  -3.41%  '
$ws.Range("D19").Value = '0.06927'
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("D20").Value = '6.586'
$ws.Range("E20").Value = '  -3.46%  '
$ws.Range("D21").Value = '17.23'
$ws.Range("E21").Value = '  -2.35%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = '12.32'
$ws.Range("E23").Value = '  -3.73%  '
$ws.Range("D24").Value = '23.324.93'
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").Value = '2.434'
$ws.Range("E25").Value = '  -2.85%  '
$ws.Range("D26").Value = '2.989'
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("D28").Value = '150.67'
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("D29").Value = '5.150'
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("D30").Value = '131.37'
$ws.Range("E30").Value = '  -2.11%  '
$ws.Range("D31").Value = '1.810.26'
$ws.Range("E31").Value = '  -1.77%  '
$ws.Range("D32").Value = '6.672'
$ws.Range("E32").Value = '  -5.29%  '
$ws.Range("D33").Value = '2.119'
$ws.Range("E33").Value = '  -4.81%  '
$ws.Range("D34").Value = '11.20'
$ws.Range("E34").Value = '  -6.95%  '
$ws.Range("D35").Value = '0.9696'
$ws.Range("E35").Value = '  -8.41%  '
$ws.Range("D36").Value = '0.02648'
$ws.Range("E36").Value = '  -5.62%  '
$ws.Range("D37").Value = '0.08728'
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").Value = '0.2409'
$ws.Range("E38").Value = '  -4.18%  '
$ws.Range("D39").Value = '5.833'
$ws.Range("E39").Value = '  -4.05%  '
$ws.Range("D40").Value = '0.06654'
$ws.Range("E40").Value = '  -5.33%  '
$ws.Range("D41").Value = '12.63'
$ws.Range("E41").Value = '  -2.77%  '
$ws.Range("D42").Value = '0.6759'
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("D43").Value = '1.292'
$ws.Range("E43").Value = '  -3.22%  '
$ws.Range("D44").Value = '15.29'
$ws.Range("E44").Value = '  -5.29%  '
$ws.Range("D45").Value = '1.002'
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = '0.6281'
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").Value = '2.225'
$ws.Range("E47").Value = '  -3.63%  '
$ws.Range("D48").Value = '3.879'
$ws.Range("E48").Value = '  -2.13%  '
$ws.Range("D49").Value = '126.45'
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("D50").Value = '0.07615'
$ws.Range("E50").Value = '  -3.96%  '
$ws.Range("D51").Value = '1.214'
$ws.Range("E51").Value = '  +1.14%  '
